$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Row 72 and row 73 effectively swap their species-find content (the
# "Svart trolldruva" find moves up to row 72, the "Rodgul
# trumpetsvamp" find moves down to row 73). Only the cells that
# actually differ are touched; shared/unchanged cells (dates,
# location info, reporter, flags, ...) are left alone. The Ost/Nord
# (Q/R) coordinates are written back as plain rounded integers
# instead of the original long decimal values.
# -----------------------------------------------------------------

# --- Row 72 gets what used to be row 73's find ---
$ws.Range("A72").Value = 112090750
$ws.Range("B72").Value = 98446
$ws.Range("E72").Value = 222771
$ws.Range("F72").Value = "Svart trolldruva"
$ws.Range("G72").Value = "Actaea spicata"
$ws.Range("H72").Value = "L."
$ws.Range("I72").Value = "'1"
$ws.Range("J72").Value = "plantor/tuvor"
$ws.Range("K72").Value = "i frukt"
$ws.Range("Q72").Value = 654798
$ws.Range("R72").Value = 6626355
$ws.Range("Z72").Value = "16:02"
$ws.Range("AB72").Value = "16:02"

# --- Row 73 gets what used to be row 72's find ---
$ws.Range("A73").Value = 112090588
$ws.Range("B73").Value = 89183
$ws.Range("E73").Value = 3215
$ws.Range("F73").Value = "Rödgul trumpetsvamp"
$ws.Range("G73").Value = "Craterellus lutescens"
$ws.Range("H73").Value = "(Fr.) Fr."
# I73/K73 become blank cells (kept, but with no content); J73 is removed entirely
$ws.Range("I73").Value = "'"
$ws.Range("J73").Value = ""
$ws.Range("K73").Value = "'"
$ws.Range("Q73").Value = 654788
$ws.Range("R73").Value = 6626334
$ws.Range("Z73").Value = "15:59"
$ws.Range("AB73").Value = "15:59"

# -----------------------------------------------------------------
# New row 74 - a newly added species observation.
# -----------------------------------------------------------------
$ws.Range("A74").Value = 112249602
$ws.Range("B74").Value = 89405
$ws.Range("C74").Value = "Ovaliderad"
$ws.Range("D74").Value = "NT"
$ws.Range("E74").Value = 1202
$ws.Range("F74").Value = "Ullticka"
$ws.Range("G74").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H74").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("I74").Value = "'"
$ws.Range("K74").Value = "'"
$ws.Range("P74").Value = "Knivsta (Knivsta), Upl"
$ws.Range("Q74").Value = 654806
$ws.Range("R74").Value = 6626334
$ws.Range("S74").Value = 5
$ws.Range("T74").Value = "Uppsala"
$ws.Range("U74").Value = "Knivsta"
$ws.Range("V74").Value = "Uppland"
$ws.Range("W74").Value = "Alsike"
$ws.Range("Y74").Value = "'2023-09-22"
$ws.Range("Z74").Value = "15:14"
$ws.Range("AA74").Value = "'2023-09-22"
$ws.Range("AB74").Value = "15:14"
$ws.Range("AD74").Value = $false
$ws.Range("AE74").Value = $false
$ws.Range("AG74").Value = $false
$ws.Range("AT74").Value = "'"
$ws.Range("AW74").Value = "Marie Kvarnström"
$ws.Range("AX74").Value = "Marie Kvarnström"
$ws.Range("AY74").Value = "'"
